$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scrape): updated price/volume
# figures, plus the WrappedEther/ShibaInu row swap (rows 17-18).
# Cells are forced to Text (NumberFormat "@") before the write so
# numeric-looking strings like '56.717.40' or '1.00' survive as literal
# text instead of being parsed into floating-point numbers, then the
# style is reset to Normal so no stray cell formatting is introduced.
$updates = @(
    @{ Cell = 'D2'; Value = '56.717.40' }
    @{ Cell = 'E2'; Value = '  -2.54%  ' }
    @{ Cell = 'D3'; Value = '2.986.57' }
    @{ Cell = 'E3'; Value = '  -4.76%  ' }
    @{ Cell = 'E4'; Value = '  -0.01%  ' }
    @{ Cell = 'D5'; Value = '498.02' }
    @{ Cell = 'E5'; Value = '  -4.97%  ' }
    @{ Cell = 'D6'; Value = '134.75' }
    @{ Cell = 'E6'; Value = '  -0.10%  ' }
    @{ Cell = 'E7'; Value = '  -0.03%  ' }
    @{ Cell = 'D8'; Value = '2.984.73' }
    @{ Cell = 'E8'; Value = '  -4.79%  ' }
    @{ Cell = 'E9'; Value = '  -4.09%  ' }
    @{ Cell = 'E10'; Value = '  +0.42%  ' }
    @{ Cell = 'E11'; Value = '  -3.69%  ' }
    @{ Cell = 'D12'; Value = '0.351' }
    @{ Cell = 'E12'; Value = '  -7.77%  ' }
    @{ Cell = 'E13'; Value = '  -0.44%  ' }
    @{ Cell = 'D14'; Value = '3.490.14' }
    @{ Cell = 'D15'; Value = '24.75' }
    @{ Cell = 'E15'; Value = '  -3.37%  ' }
    @{ Cell = 'D16'; Value = '56.619.42' }
    @{ Cell = 'E16'; Value = '  -2.68%  ' }
    @{ Cell = 'B17'; Value = 'ShibaInu' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Cell = 'D17'; Value = '0.0000147' }
    @{ Cell = 'E17'; Value = '  -3.45%  ' }
    @{ Cell = 'B18'; Value = 'WrappedEther' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D18'; Value = '2.979.62' }
    @{ Cell = 'E18'; Value = '  -5.05%  ' }
    @{ Cell = 'D19'; Value = '5.81' }
    @{ Cell = 'E19'; Value = '  +0.24%  ' }
    @{ Cell = 'D20'; Value = '12.33' }
    @{ Cell = 'E20'; Value = '  -5.71%  ' }
    @{ Cell = 'D21'; Value = '7.76' }
    @{ Cell = 'E21'; Value = '  -2.28%  ' }
    @{ Cell = 'D22'; Value = '325.82' }
    @{ Cell = 'E22'; Value = '  -5.34%  ' }
    @{ Cell = 'D23'; Value = '0.998' }
    @{ Cell = 'E23'; Value = '  -0.29%  ' }
    @{ Cell = 'D24'; Value = '0.465' }
    @{ Cell = 'E24'; Value = '  -8.29%  ' }
    @{ Cell = 'D25'; Value = '61.31' }
    @{ Cell = 'E25'; Value = '  -10.64%  ' }
    @{ Cell = 'D26'; Value = '0.995' }
    @{ Cell = 'E26'; Value = '  -0.61%  ' }
    @{ Cell = 'E27'; Value = '  -3.39%  ' }
    @{ Cell = 'D28'; Value = '0.0₃0910' }
    @{ Cell = 'E28'; Value = '  -5.01%  ' }
    @{ Cell = 'E29'; Value = '  -0.01%  ' }
    @{ Cell = 'D30'; Value = '6.51' }
    @{ Cell = 'E30'; Value = '  -4.37%  ' }
    @{ Cell = 'D31'; Value = '6.79' }
    @{ Cell = 'E31'; Value = '  -1.03%  ' }
    @{ Cell = 'D32'; Value = '1.18' }
    @{ Cell = 'E32'; Value = '  -4.03%  ' }
    @{ Cell = 'E33'; Value = '  -6.65%  ' }
    @{ Cell = 'D34'; Value = '19.98' }
    @{ Cell = 'E34'; Value = '  -7.04%  ' }
    @{ Cell = 'D35'; Value = '154.58' }
    @{ Cell = 'E35'; Value = '  -1.46%  ' }
    @{ Cell = 'D36'; Value = '4.48' }
    @{ Cell = 'E36'; Value = '  -6.78%  ' }
    @{ Cell = 'E37'; Value = '  -6.75%  ' }
    @{ Cell = 'D38'; Value = '5.62' }
    @{ Cell = 'E38'; Value = '  -9.81%  ' }
    @{ Cell = 'D39'; Value = '0.0674' }
    @{ Cell = 'E39'; Value = '  -2.42%  ' }
    @{ Cell = 'D40'; Value = '23.39' }
    @{ Cell = 'E40'; Value = '  -3.88%  ' }
    @{ Cell = 'D41'; Value = '3.015.37' }
    @{ Cell = 'E41'; Value = '  -4.84%  ' }
    @{ Cell = 'D42'; Value = '36.69' }
    @{ Cell = 'E42'; Value = '  -9.20%  ' }
    @{ Cell = 'D43'; Value = '1.00' }
    @{ Cell = 'E43'; Value = '  +0.02%  ' }
    @{ Cell = 'E44'; Value = '  -6.49%  ' }
    @{ Cell = 'E45'; Value = '  -7.95%  ' }
    @{ Cell = 'D46'; Value = '1.41' }
    @{ Cell = 'E46'; Value = '  -2.08%  ' }
    @{ Cell = 'D47'; Value = '2.205.13' }
    @{ Cell = 'E47'; Value = '  -2.34%  ' }
    @{ Cell = 'D48'; Value = '3.56' }
    @{ Cell = 'E48'; Value = '  -8.85%  ' }
    @{ Cell = 'D49'; Value = '1.96' }
    @{ Cell = 'E49'; Value = '  +6.47%  ' }
    @{ Cell = 'D50'; Value = '0.0237' }
    @{ Cell = 'E50'; Value = '  +1.92%  ' }
    @{ Cell = 'D51'; Value = '5.72' }
    @{ Cell = 'E51'; Value = '  -7.51%  ' }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = '@'
    $range.Value = $update.Value
    $range.Style = 'Normal'
}
